# Update births (column E) values per the data push
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{Row=2; Value=17},
    @{Row=9; Value=29},
    @{Row=10; Value=20},
    @{Row=11; Value=28},
    @{Row=12; Value=21},
    @{Row=13; Value=28},
    @{Row=16; Value=26},
    @{Row=17; Value=26},
    @{Row=18; Value=32},
    @{Row=19; Value=28},
    @{Row=22; Value=31},
    @{Row=23; Value=19},
    @{Row=27; Value=22},
    @{Row=28; Value=23},
    @{Row=32; Value=24},
    @{Row=35; Value=23},
    @{Row=38; Value=20},
    @{Row=49; Value=25},
    @{Row=51; Value=24},
    @{Row=53; Value=24},
    @{Row=63; Value=22},
    @{Row=66; Value=29},
    @{Row=68; Value=26},
    @{Row=69; Value=29},
    @{Row=71; Value=23},
    @{Row=73; Value=32},
    @{Row=74; Value=23},
    @{Row=76; Value=30},
    @{Row=82; Value=23},
    @{Row=84; Value=21},
    @{Row=85; Value=25},
    @{Row=88; Value=31},
    @{Row=89; Value=19},
    @{Row=92; Value=18},
    @{Row=94; Value=26},
    @{Row=96; Value=29},
    @{Row=97; Value=20},
    @{Row=101; Value=22},
    @{Row=103; Value=26},
    @{Row=104; Value=31},
    @{Row=109; Value=16},
    @{Row=110; Value=26},
    @{Row=111; Value=26},
    @{Row=117; Value=28},
    @{Row=118; Value=27},
    @{Row=120; Value=26},
    @{Row=123; Value=18},
    @{Row=128; Value=20},
    @{Row=131; Value=31},
    @{Row=133; Value=24},
    @{Row=134; Value=28},
    @{Row=135; Value=31},
    @{Row=136; Value=23},
    @{Row=139; Value=19},
    @{Row=145; Value=20},
    @{Row=147; Value=33},
    @{Row=148; Value=30},
    @{Row=150; Value=22},
    @{Row=151; Value=28},
    @{Row=157; Value=21},
    @{Row=165; Value=20},
    @{Row=167; Value=21},
    @{Row=171; Value=24},
    @{Row=172; Value=25},
    @{Row=173; Value=27},
    @{Row=177; Value=18},
    @{Row=182; Value=27},
    @{Row=183; Value=26},
    @{Row=184; Value=21},
    @{Row=186; Value=31},
    @{Row=187; Value=23},
    @{Row=188; Value=21},
    @{Row=189; Value=22},
    @{Row=190; Value=31},
    @{Row=192; Value=21},
    @{Row=194; Value=22},
    @{Row=195; Value=28},
    @{Row=196; Value=28},
    @{Row=197; Value=24},
    @{Row=199; Value=17},
    @{Row=201; Value=18},
    @{Row=209; Value=24},
    @{Row=210; Value=25},
    @{Row=211; Value=16},
    @{Row=213; Value=28},
    @{Row=215; Value=25},
    @{Row=216; Value=26},
    @{Row=217; Value=25},
    @{Row=218; Value=44},
    @{Row=223; Value=40},
    @{Row=225; Value=26},
    @{Row=229; Value=38},
    @{Row=230; Value=30},
    @{Row=232; Value=35},
    @{Row=236; Value=33},
    @{Row=237; Value=32},
    @{Row=239; Value=35},
    @{Row=242; Value=24},
    @{Row=244; Value=37},
    @{Row=248; Value=30},
    @{Row=249; Value=30},
    @{Row=250; Value=33},
    @{Row=251; Value=23},
    @{Row=252; Value=25},
    @{Row=254; Value=32},
    @{Row=258; Value=32},
    @{Row=260; Value=28},
    @{Row=262; Value=28},
    @{Row=266; Value=28},
    @{Row=267; Value=34},
    @{Row=269; Value=38},
    @{Row=270; Value=30},
    @{Row=273; Value=33},
    @{Row=274; Value=31},
    @{Row=275; Value=24},
    @{Row=276; Value=32},
    @{Row=279; Value=45},
    @{Row=281; Value=31},
    @{Row=283; Value=32},
    @{Row=285; Value=28},
    @{Row=286; Value=26},
    @{Row=289; Value=25},
    @{Row=292; Value=30},
    @{Row=293; Value=27},
    @{Row=295; Value=21},
    @{Row=301; Value=28},
    @{Row=304; Value=24},
    @{Row=307; Value=23},
    @{Row=311; Value=18},
    @{Row=320; Value=16},
    @{Row=322; Value=23},
    @{Row=324; Value=42},
    @{Row=325; Value=25},
    @{Row=329; Value=37},
    @{Row=334; Value=26},
    @{Row=340; Value=21},
    @{Row=343; Value=26},
    @{Row=345; Value=31},
    @{Row=346; Value=25},
    @{Row=351; Value=26},
    @{Row=352; Value=27},
    @{Row=353; Value=27},
    @{Row=354; Value=28},
    @{Row=357; Value=35},
    @{Row=358; Value=35},
    @{Row=361; Value=26},
    @{Row=364; Value=17}
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 5).Value = $u.Value
}

$wb.Save()
